$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="23.763.66"'
$ws.Range("E2").Formula = '="  +1.78%  "'
$ws.Range("D3").Formula = '="1.655.85"'
$ws.Range("E3").Formula = '="  +1.87%  "'
$ws.Range("D4").Formula = '="1.002"'
$ws.Range("E4").Formula = '="  +0.20%  "'
$ws.Range("E5").Formula = '="  +0.08%  "'
$ws.Range("D6").Formula = '="304.17"'
$ws.Range("E6").Formula = '="  +0.52%  "'
$ws.Range("D7").Formula = '="0.3817"'
$ws.Range("E7").Formula = '="  +1.49%  "'
$ws.Range("E8").Formula = '="  -0.33%  "'
$ws.Range("D9").Formula = '="51.30"'
$ws.Range("E9").Formula = '="  -0.27%  "'
$ws.Range("D10").Formula = '="1.256"'
$ws.Range("E10").Formula = '="  +2.80%  "'
$ws.Range("D11").Formula = '="0.08221"'
$ws.Range("E11").Formula = '="  +1.02%  "'
$ws.Range("E12").Formula = '="  +0.15%  "'
$ws.Range("D13").Formula = '="22.69"'
$ws.Range("E13").Formula = '="  +1.83%  "'
$ws.Range("D14").Formula = '="6.542"'
$ws.Range("E14").Formula = '="  +1.18%  "'
$ws.Range("D15").Formula = '="7.463"'
$ws.Range("E15").Formula = '="  +2.43%  "'
$ws.Range("D16").Formula = '="0.00001240"'
$ws.Range("E16").Formula = '="  +0.14%  "'
$ws.Range("D17").Formula = '="1.637.52"'
$ws.Range("E17").Formula = '="  +0.90%  "'
$ws.Range("D18").Formula = '="97.86"'
$ws.Range("E18").Formula = '="  +3.90%  "'
$ws.Range("D19").Formula = '="0.06984"'
$ws.Range("E19").Formula = '="  +0.91%  "'
$ws.Range("D20").Formula = '="6.779"'
$ws.Range("E20").Formula = '="  +3.56%  "'
$ws.Range("D21").Formula = '="17.77"'
$ws.Range("E21").Formula = '="  +1.28%  "'
$ws.Range("D22").Formula = '="1.000"'
$ws.Range("E22").Formula = '="  -0.02%  "'
$ws.Range("D23").Formula = '="12.76"'
$ws.Range("E23").Formula = '="  +2.03%  "'
$ws.Range("D24").Formula = '="2.573"'
$ws.Range("E24").Formula = '="  +3.26%  "'
$ws.Range("D25").Formula = '="23.768.87"'
$ws.Range("E25").Formula = '="  +1.86%  "'
$ws.Range("D26").Formula = '="3.081"'
$ws.Range("E26").Formula = '="  -0.14%  "'
$ws.Range("D27").Formula = '="21.30"'
$ws.Range("E27").Formula = '="  +0.83%  "'
$ws.Range("D28").Formula = '="150.91"'
$ws.Range("E28").Formula = '="  +0.69%  "'
$ws.Range("D29").Formula = '="5.225"'
$ws.Range("E29").Formula = '="  -0.94%  "'
$ws.Range("D30").Formula = '="134.33"'
$ws.Range("E30").Formula = '="  +1.21%  "'
$ws.Range("D31").Formula = '="1.822.01"'
$ws.Range("E31").Formula = '="  +1.11%  "'
$ws.Range("D32").Formula = '="6.947"'
$ws.Range("E32").Formula = '="  +3.37%  "'
$ws.Range("D33").Formula = '="2.180"'
$ws.Range("E33").Formula = '="  +2.62%  "'
$ws.Range("E34").Formula = '="  +1.93%  "'
$ws.Range("D35").Formula = '="11.86"'
$ws.Range("E35").Formula = '="  +6.35%  "'
$ws.Range("E36").Formula = '="  +2.83%  "'
$ws.Range("D37").Formula = '="6.160"'
$ws.Range("E37").Formula = '="  +2.73%  "'
$ws.Range("D38").Formula = '="0.2521"'
$ws.Range("E38").Formula = '="  +1.56%  "'
$ws.Range("D39").Formula = '="0.08830"'
$ws.Range("E39").Formula = '="  +0.76%  "'
$ws.Range("D40").Formula = '="0.07180"'
$ws.Range("E40").Formula = '="  +1.09%  "'
$ws.Range("D41").Formula = '="13.27"'
$ws.Range("E41").Formula = '="  +9.98%  "'
$ws.Range("D42").Formula = '="0.7073"'
$ws.Range("E42").Formula = '="  +1.15%  "'
$ws.Range("E43").Formula = '="  +0.31%  "'
$ws.Range("D44").Formula = '="16.08"'
$ws.Range("E44").Formula = '="  +1.51%  "'
$ws.Range("D45").Formula = '="0.6553"'
$ws.Range("E45").Formula = '="  +1.57%  "'
$ws.Range("E46").Formula = '="  +2.72%  "'
$ws.Range("D47").Formula = '="1.000"'
$ws.Range("E47").Formula = '="  +0.06%  "'
$ws.Range("D48").Formula = '="3.960"'
$ws.Range("E48").Formula = '="  +0.00%  "'
$ws.Range("D49").Formula = '="0.07973"'
$ws.Range("E49").Formula = '="  +0.02%  "'
$ws.Range("D50").Formula = '="128.62"'
$ws.Range("E50").Formula = '="  +1.93%  "'
$ws.Range("D51").Formula = '="1.193"'
$ws.Range("E51").Formula = '="  +0.56%  "'

$rng = $ws.Range("D2:E51")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0
